$wb = $excel.ActiveWorkbook

# --- Workbook-level window position (xWindow 0 -> -20) ---
$wb.Windows.Item(1).Left = -20

# --- Sheets ---
$ws1 = $wb.Worksheets.Item("D7SCFA")
$ws3 = $wb.Worksheets.Item("BioSample_ID")

# --- D7SCFA: column J (Bacteroides group) relabel "LB" -> "NB" ---
# Rows whose J value was the shared string "LB" (67) become the new
# label "NB"; rows already labelled "HB" are left untouched.
$nbRows = @(2,3,4,6,8,11,14,15,16,19,20,24,25,26,27,28,34)
foreach ($r in $nbRows) {
    $ws1.Cells.Item($r, 10).Value2 = "NB"
}

# --- Row 13 lost its one-off highlight formatting (fill) that made it
#     stand out from the rest of the table; bring it back in line with
#     the surrounding rows, then restore the per-column formats that are
#     still meant to be there (I: category font, J: Bacteroides font). ---
$ws1.Rows.Item(13).ClearFormats()

$ws1.Range("I14").Copy()
$ws1.Range("I13").PasteSpecial(-4122)

$ws1.Range("J14").Copy()
$ws1.Range("J13").PasteSpecial(-4122)
$ws1.Range("J13").Value2 = "NB"

$ws1.Application.CutCopyMode = $false

# --- Selections (cursor position saved with the sheet) ---
$ws3.Range("Q8").Select()
$ws1.Activate()
$ws1.Range("J35").Select()
